# "Generate Report for Handoff" -- the localization-status report is
# regenerated: the workbook-wide status moves from "In Translation" to
# "Ready for handoff", the handoff timestamps advance a few seconds, and
# the Status/Latest Handoff Datetime columns widen to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$overview.Range("F2").Value = "Ready for handoff"   # de-de status column
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status
$dede.Range("C2").Value     = "Ready for handoff"   # Status

# --- Refreshed handoff timestamps --------------------------------------
$overview.Range("G2").Value = "2016-08-23 00:55:33" # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-08-23 00:55:33" # Latest Handoff Datetime

$zhcn.Range("H2").Value     = "2016-08-23 00:55:29" # Latest Handoff Datetime

# --- Widen the Status / Latest Handoff Datetime columns to fit the new
#     longer text. The COM ColumnWidth setter snaps to the nearest 1/6
#     of a character, so use the value that lands closest to the wider
#     target width used by the regenerated report.
$overview.Columns.Item(5).ColumnWidth = 16.3   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 16.3   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
$dede.Columns.Item(3).ColumnWidth     = 16.3   # column C (Status)
